# Daily attendance processing - rotate the "Recorded By" (column G) list so
# that the first recorder in the comma-separated list is moved to the end.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 157 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -eq $null) { continue }
    $text = [string]$val
    if ($text -eq "") { continue }

    $parts = $text -split ", "
    if ($parts.Length -gt 1) {
        $rotated = $parts[1..($parts.Length - 1)] + $parts[0]
        $newText = $rotated -join ", "
        if ($newText -ne $text) {
            $cell.Value2 = $newText
        }
    }
}
